# regen sval data to filter save games
# Updates the numeric stat columns (B:E, G) for rows 2-6 on Sheet1 with
# newly-regenerated values (F / "Win" flags and the date labels in column A
# are unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    3 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    4 = @(0.3048080303191223, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 1.919867272924993)
    5 = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 11.80884686099532)
    6 = @(3.230985683306322, 1.667794583268128, 26.21740644021617, 8.660232485948974, 39.7764191927396)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B: TB
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C: d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D: K
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E: IP
    $ws.Cells.Item($row, 7).Value = $vals[4]  # G: sum
}
